$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet (and workbook view title) from "Through 2021-09-29" to "Through 2021-09-30"
$ws.Name = "Through 2021-09-30"

# Row 11: "September (through 09-29)" -> "September (through 09-30)" plus updated counts/rates
$ws.Range("A11").Value = "September (through 09-30)"

$ws.Range("C11").Value = 31
$ws.Range("D11").Value = 0.0312

$ws.Range("F11").Value = 43
$ws.Range("G11").Value = 0.1042

$ws.Range("I11").Value = 71
$ws.Range("J11").Value = 0.0658

$ws.Range("L11").Value = 54
$ws.Range("M11").Value = 0.069

$ws.Range("O11").Value = 66
$ws.Range("P11").Value = 0.0959

$ws.Range("R11").Value = 112
$ws.Range("S11").Value = 0.0427

$ws.Range("U11").Value = 176
$ws.Range("V11").Value = 0.0112

# Row 12 (Total): updated counts/rates
$ws.Range("C12").Value = 196
$ws.Range("D12").Value = 0.1327

$ws.Range("F12").Value = 383
$ws.Range("G12").Value = 0.1072

$ws.Range("I12").Value = 577
$ws.Range("J12").Value = 0.0797

$ws.Range("L12").Value = 487
$ws.Range("M12").Value = 0.1113

$ws.Range("O12").Value = 379
$ws.Range("P12").Value = 0.1019

$ws.Range("R12").Value = 848
$ws.Range("S12").Value = 0.0588

$ws.Range("U12").Value = 1170
$ws.Range("V12").Value = 0.0625

$wb.Save()
